$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 9; existing rows 9-27 shift down to 10-28,
# preserving every column's contents (A,B,C,E,F,G,H,I,L,N,O,Q,R stay the
# same "Achicoria / Vega Central Mapocho" template; D,J,K,M,P are the
# per-row observations that ride along with their original rows).
$ws.Rows(9).Insert()

# Populate the freshly inserted row 9 with the new weekly observation.
$ws.Cells.Item(9, 1).Value = 9
$ws.Cells.Item(9, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(9, 3).Value = "Metropolitana"
$ws.Cells.Item(9, 4).Value = [DateTime]"2022-06-03"
$ws.Cells.Item(9, 5).Value = 13
$ws.Cells.Item(9, 6).Value = 100112010
$ws.Cells.Item(9, 7).Value = "Achicoria"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 70
$ws.Cells.Item(9, 11).Value = 5000
$ws.Cells.Item(9, 12).Value = 6000
$ws.Cells.Item(9, 13).Value = 5500
$ws.Cells.Item(9, 14).Value = "`$/caja 16 unidades"
$ws.Cells.Item(9, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(9, 16).Value = 344
$ws.Cells.Item(9, 17).Value = 16
$ws.Cells.Item(9, 18).Value = "Hortaliza"
